# Hortaliza, Feria Lagunitas de Puerto Montt - Papa
# A new weekly price-survey row is inserted at row 366 (pushing the
# existing rows 366-411 down to 367-412), and the final row (old row 411,
# now sitting at row 412) is duplicated onto a brand-new row 413.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank row at 366; rows 366-411 shift down to 367-412.
$ws.Rows.Item(366).Insert()

# 2) Populate the newly inserted row 366 with the new observation.
$ws.Cells.Item(366, 1).Value = 4
$ws.Cells.Item(366, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(366, 3).Value = "Los Lagos"
$ws.Cells.Item(366, 4).Value = 44748
$ws.Cells.Item(366, 5).Value = 10
$ws.Cells.Item(366, 6).Value = 100114001
$ws.Cells.Item(366, 7).Value = "Papa"
$ws.Cells.Item(366, 8).Value = "Patagonia"
$ws.Cells.Item(366, 9).Value = "1a (guarda)"
$ws.Cells.Item(366, 10).Value = 150
$ws.Cells.Item(366, 11).Value = 7000
$ws.Cells.Item(366, 12).Value = 8000
$ws.Cells.Item(366, 13).Value = 7467
$ws.Cells.Item(366, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(366, 15).Value = "Provincia de Llanquihue"
$ws.Cells.Item(366, 16).Value = 299
$ws.Cells.Item(366, 17).Value = 25
$ws.Cells.Item(366, 18).Value = "Hortaliza"

# Match the date-formatted style already used by the other D-column cells.
$ws.Cells.Item(366, 4).NumberFormat = $ws.Cells.Item(367, 4).NumberFormat

# 3) Append a new row 413 duplicating row 412 (the old last row, 411,
#    which is now at 412 after the insert above).
$ws.Cells.Item(413, 1).Value = 4
$ws.Cells.Item(413, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(413, 3).Value = "Los Lagos"
$ws.Cells.Item(413, 4).Value = 44194
$ws.Cells.Item(413, 5).Value = 10
$ws.Cells.Item(413, 6).Value = 100114001
$ws.Cells.Item(413, 7).Value = "Papa"
$ws.Cells.Item(413, 8).Value = "Pukará"
$ws.Cells.Item(413, 9).Value = "1a nueva(o)"
$ws.Cells.Item(413, 10).Value = 600
$ws.Cells.Item(413, 11).Value = 14000
$ws.Cells.Item(413, 12).Value = 15000
$ws.Cells.Item(413, 13).Value = 14500
$ws.Cells.Item(413, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(413, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(413, 16).Value = 580
$ws.Cells.Item(413, 17).Value = 25
$ws.Cells.Item(413, 18).Value = "Hortaliza"

$ws.Cells.Item(413, 4).NumberFormat = $ws.Cells.Item(412, 4).NumberFormat
